# Update the acquisition timestamp (取得日時) in column A for rows 2-13
# on the "ランサーズ" sheet from "2025-10-22 01:21:01" to "2025-10-22 01:51:12".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq "2025-10-22 01:21:01") {
        $cell.Value = "2025-10-22 01:51:12"
    }
}
